$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old layout was: A=mapsfrom, B=meaning (free text), C=mapsto (numeric code).
# Drop the "meaning" column entirely; this shifts the old "mapsto" column (C)
# left into B, carrying over its original formatting.
$ws.Columns("B").Delete()

# Column B (formerly C) now holds the numeric mapsto codes. Replace them with
# clearer/shorter string labels, and collapse the old S3/S4 "deep sleep"
# stages into a single "sws" (slow-wave sleep) bucket.
$ws.Range("B1").Value = "mapsto"
$ws.Range("B2").Value = "wake"
$ws.Range("B3").Value = "stage1"
$ws.Range("B4").Value = "stage2"
$ws.Range("B5").Value = "sws"
$ws.Range("B6").Value = "sws"
$ws.Range("B7").Value = "rem"
$ws.Range("B8").Value = "rem"

# Match the refreshed view/selection from the edit.
$ws.Range("B9").Select()
